$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3194513333333333
$ws.Range("H2").Value = 0.9583539999999999
$ws.Range("I2").Value = 0.01095865642710367
$ws.Range("J2").Value = 0.01095865642710367
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 3.560874817161334
$ws.Range("R2").Value = 32.047873354452
$ws.Range("S2").Value = 0.002843558838783867
$ws.Range("T2").Value = 0.002843558838783867

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3194513333333333
$ws.Range("H3").Value = 0.9583539999999999
$ws.Range("I3").Value = 0.01095865642710367
$ws.Range("J3").Value = 0.01095865642710367
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 8.844168930645999
$ws.Range("R3").Value = 79.59752037581399
$ws.Range("S3").Value = 0.007062566370834783
$ws.Range("T3").Value = 0.007062566370834784

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3194513333333333
$ws.Range("H4").Value = 0.9583539999999999
$ws.Range("I4").Value = 0.01095865642710367
$ws.Range("J4").Value = 0.01095865642710367
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 1.318042677893556
$ws.Range("R4").Value = 11.862384101042
$ws.Range("S4").Value = 0.001052531217485023
$ws.Range("T4").Value = 0.001052531217485023

$ws.Range("I5").Value = 0.9713235907985359
$ws.Range("J5").Value = 0.971323590798536
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 315.6191396998987
$ws.Range("R5").Value = 2840.572257299088
$ws.Range("S5").Value = 0.2520396364560952
$ws.Range("T5").Value = 0.2520396364560952

$ws.Range("I6").Value = 0.9713235907985359
$ws.Range("J6").Value = 0.971323590798536
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.6259925542154537
$ws.Range("T6").Value = 0.6259925542154537

$ws.Range("I7").Value = 0.9713235907985359
$ws.Range("J7").Value = 0.971323590798536
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.09329140012698708
$ws.Range("T7").Value = 0.09329140012698708

$ws.Range("I8").Value = 0.01771775277436037
$ws.Range("J8").Value = 0.01771775277436037
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 5.757156462618001
$ws.Range("R8").Value = 51.81440816356201
$ws.Range("S8").Value = 0.004597413272334466
$ws.Range("T8").Value = 0.004597413272334466

$ws.Range("I9").Value = 0.01771775277436037
$ws.Range("J9").Value = 0.01771775277436037
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.01141862652080921
$ws.Range("T9").Value = 0.01141862652080921

$ws.Range("I10").Value = 0.01771775277436037
$ws.Range("J10").Value = 0.01771775277436037
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.001701712981216703
$ws.Range("T10").Value = 0.001701712981216703

